$wb = $excel.ActiveWorkbook

# --- mon ---
$ws = $wb.Worksheets.Item("mon")
$ws.Range("G2").Value = "EDS121"
$ws.Range("I2").Value = "PHY121"
$ws.Range("J2").Value = "PHY121"
$ws.Range("I3").Value = "CSC425"
$ws.Range("J3").Value = "CSC425"
$ws.Range("F6").Value = "CSC225"
$ws.Range("C7").Value = "MAT121"
$ws.Range("D7").Value = "MAT121"
$ws.Range("E10").Value = "MAT122"
$ws.Range("F10").Value = "MAT122"
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = ""
$ws.Range("J15").Value = ""
$ws.Range("K15").Value = ""
$ws.Range("C18").Value = "PHY122"
$ws.Range("D18").Value = "PHY122"
$ws.Range("K19").Value = "MAT225"
$ws.Range("E20").Value = ""
$ws.Range("E21").Value = "TMC421"
$ws.Range("F21").Value = ""
$ws.Range("H21").Value = ""
$ws.Range("I21").Value = ""
$ws.Range("J25").Value = ""
$ws.Range("H27").Value = "EDS421"

# --- tue ---
$ws = $wb.Worksheets.Item("tue")
$ws.Range("C2").Value = "CIT224"
$ws.Range("D2").Value = "CIT224"
$ws.Range("D3").Value = "CSC221"
$ws.Range("J6").Value = "CSC443"
$ws.Range("K6").Value = "CSC443"
$ws.Range("J8").Value = ""
$ws.Range("K8").Value = ""
$ws.Range("G13").Value = ""
$ws.Range("C15").Value = "CSC446"
$ws.Range("D15").Value = "CSC446"
$ws.Range("I16").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("F17").Value = ""
$ws.Range("J17").Value = "CSC444"
$ws.Range("K17").Value = "CSC444"
$ws.Range("D18").Value = "MAT121"
$ws.Range("J19").Value = "CIT221"
$ws.Range("K19").Value = "CIT221"
$ws.Range("G21").Value = ""
$ws.Range("H21").Value = ""
$ws.Range("J21").Value = "CSC223"
$ws.Range("K21").Value = "CSC223"
$ws.Range("G24").Value = "CSC227"
$ws.Range("H24").Value = "CSC227"

# --- wed ---
$ws = $wb.Worksheets.Item("wed")
$ws.Range("I4").Value = "CSC425"
$ws.Range("E8").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F13").Value = "CIT141"
$ws.Range("G13").Value = "CIT141"
$ws.Range("J14").Value = "MAT226"
$ws.Range("K14").Value = "MAT226"
$ws.Range("F16").Value = ""
$ws.Range("G16").Value = ""
$ws.Range("C20").Value = ""
$ws.Range("D20").Value = ""
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = ""
$ws.Range("G20").Value = ""
$ws.Range("H20").Value = ""
$ws.Range("J21").Value = ""
$ws.Range("K21").Value = ""
$ws.Range("E25").Value = "TMC121"
$ws.Range("B26").Value = "CSC221"
$ws.Range("C26").Value = "CSC221"
$ws.Range("D28").Value = "PHY129"
$ws.Range("E28").Value = ""
$ws.Range("F28").Value = ""

# --- thur ---
$ws = $wb.Worksheets.Item("thur")
$ws.Range("H4").Value = "CST121"
$ws.Range("I4").Value = "CST121"
$ws.Range("F7").Value = "MIS221"
$ws.Range("G8").Value = "CSC442"
$ws.Range("H8").Value = "CSC442"
$ws.Range("J14").Value = "MAT226"
$ws.Range("G16").Value = "CIT121"
$ws.Range("H16").Value = "CIT121"
$ws.Range("E19").Value = "CSC224"
$ws.Range("F19").Value = "CSC224"
$ws.Range("E20").Value = "CSC223"
$ws.Range("I20").Value = ""
$ws.Range("J20").Value = ""
$ws.Range("J23").Value = "MAT229"
$ws.Range("K23").Value = "MAT229"
$ws.Range("G26").Value = "CSC423"
$ws.Range("H26").Value = "CSC423"
$ws.Range("G27").Value = "GST122"
$ws.Range("H27").Value = "GST122"
$ws.Range("F28").Value = "TMC221"
$ws.Range("H29").Value = "PHY121"

# --- fri ---
$ws = $wb.Worksheets.Item("fri")
$ws.Range("D2").Value = "CSC226"
$ws.Range("E2").Value = "CSC226"
$ws.Range("D4").Value = ""
$ws.Range("F7").Value = "CSC125"
$ws.Range("G7").Value = "CSC125"
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("F19").Value = "CSC441"
$ws.Range("G19").Value = "CSC441"
$ws.Range("F20").Value = "CSC121"
$ws.Range("G20").Value = "CSC121"
$ws.Range("B26").Value = "CSC423"
$ws.Range("E28").Value = "MAT225"
$ws.Range("F28").Value = "MAT225"
